$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 132.9375
$ws.Range("I9").Value = 78.75
$ws.Range("J9").Value = 295.5
$ws.Range("K9").Value = 78.75
$ws.Range("L9").Value = 295.5
$ws.Range("M9").Value = 90.25
$ws.Range("N9").Value = -633.5

# Row 12
$ws.Range("H12").Value = 95.23077000000001
$ws.Range("I12").Value = 94.833336
$ws.Range("K12").Value = 94.833336
$ws.Range("M12").Value = 75.166664

# Row 29
$ws.Range("H29").Value = 350.5
$ws.Range("I29").Value = 100.75
$ws.Range("J29").Value = 850
$ws.Range("K29").Value = 302.25
$ws.Range("L29").Value = 2550
$ws.Range("M29").Value = -21.25
$ws.Range("N29").Value = -3112

# Row 38
$ws.Range("H38").Value = 2043.2
$ws.Range("I38").Value = 137.25
$ws.Range("J38").Value = 4221.4287
$ws.Range("K38").Value = 411.75
$ws.Range("L38").Value = 12664.2861
$ws.Range("M38").Value = -39.75
$ws.Range("N38").Value = -13408.2861

# Row 58
$ws.Range("H58").Value = 5352.857
$ws.Range("I58").Value = 80
$ws.Range("J58").Value = 6790.909
$ws.Range("K58").Value = 240
$ws.Range("L58").Value = 20372.727
$ws.Range("M58").Value = -90
$ws.Range("N58").Value = -20672.727

# Row 87
$ws.Range("H87").Value = 11949.58
$ws.Range("J87").Value = 11949.58
$ws.Range("L87").Value = 11949.58
$ws.Range("N87").Value = -14445.58

# Row 90
$ws.Range("H90").Value = 11949.58
$ws.Range("J90").Value = 11949.58
$ws.Range("L90").Value = 35848.74
$ws.Range("N90").Value = -48328.74

# Row 112
$ws.Range("H112").Value = 1201.0416
$ws.Range("I112").Value = 613.3333
$ws.Range("J112").Value = 1788.75
$ws.Range("K112").Value = 1839.9999
$ws.Range("L112").Value = 5366.25
$ws.Range("M112").Value = -731.9999
$ws.Range("N112").Value = -7582.25

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1112466.5
$ws.Range("I45").Value = 2501049.8
$ws.Range("J45").Value = 1600
$ws.Range("K45").Value = 2501049.8
$ws.Range("L45").Value = 1600
$ws.Range("M45").Value = -2500672.8
$ws.Range("N45").Value = -2354

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 847.5
$ws.Range("I107").Value = 777
$ws.Range("K107").Value = 777
$ws.Range("M107").Value = 1143

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 6635.931
$ws.Range("I86").Value = 8455.764999999999
$ws.Range("J86").Value = 4057.8333
$ws.Range("K86").Value = 8455.764999999999
$ws.Range("L86").Value = 4057.8333
$ws.Range("M86").Value = -7332.764999999999
$ws.Range("N86").Value = -6303.8333

# Row 89
$ws.Range("H89").Value = 6635.931
$ws.Range("I89").Value = 8455.764999999999
$ws.Range("J89").Value = 4057.8333
$ws.Range("K89").Value = 42278.825
$ws.Range("L89").Value = 20289.1665
$ws.Range("M89").Value = -36662.825
$ws.Range("N89").Value = -31521.1665

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 3812.5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 3812.5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 11437.5
$ws.Range("N9").Value = -11885.5
$ws.Range("M9").ClearContents()

# Row 20
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3454

# Row 21
$ws.Range("H21").Value = 1400
$ws.Range("J21").Value = 1400
$ws.Range("L21").Value = 4200
$ws.Range("N21").Value = -4546

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 2000
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2490

# Row 43
$ws.Range("H43").Value = 710
$ws.Range("I43").Value = 710
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 710
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -559
$ws.Range("N43").ClearContents()

# Row 46
$ws.Range("H46").Value = 19946
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 19946
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 19946
$ws.Range("N46").Value = -20258
$ws.Range("M46").ClearContents()

# Row 57
$ws.Range("H57").Value = 34511
$ws.Range("J57").Value = 34511
$ws.Range("L57").Value = 34511
$ws.Range("N57").Value = -36151

# Row 64
$ws.Range("H64").Value = 32000
$ws.Range("J64").Value = 32000
$ws.Range("L64").Value = 32000
$ws.Range("N64").Value = -32496

# Row 67
$ws.Range("H67").Value = 32000
$ws.Range("J67").Value = 32000
$ws.Range("L67").Value = 32000
$ws.Range("N67").Value = -33716

# Row 80
$ws.Range("H80").Value = 16362.875
$ws.Range("I80").Value = 6000
$ws.Range("J80").Value = 29686.572
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 29686.572
$ws.Range("M80").Value = -5002
$ws.Range("N80").Value = -31682.572

# Row 83
$ws.Range("H83").Value = 16362.875
$ws.Range("I83").Value = 6000
$ws.Range("J83").Value = 29686.572
$ws.Range("K83").Value = 30000
$ws.Range("L83").Value = 148432.86
$ws.Range("M83").Value = -25008
$ws.Range("N83").Value = -158416.86

# Row 102
$ws.Range("H102").Value = 17384.834
$ws.Range("I102").Value = 20601.8
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 20601.8
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = -18979.8
$ws.Range("N102").Value = -4544

# Row 126
$ws.Range("H126").Value = 17385.334
$ws.Range("I126").Value = 17385.334
$ws.Range("K126").Value = 52156.00199999999
$ws.Range("M126").Value = -49686.00199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 5000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452

# Row 22
$ws.Range("H22").Value = 2691.3125
$ws.Range("I22").Value = 1921.4286
$ws.Range("J22").Value = 3290.111
$ws.Range("K22").Value = 1921.4286
$ws.Range("L22").Value = 3290.111
$ws.Range("M22").Value = -1626.4286
$ws.Range("N22").Value = -3880.111

# Row 27
$ws.Range("H27").Value = 2691.3125
$ws.Range("I27").Value = 1921.4286
$ws.Range("J27").Value = 3290.111
$ws.Range("K27").Value = 1921.4286
$ws.Range("L27").Value = 3290.111
$ws.Range("M27").Value = -1814.4286
$ws.Range("N27").Value = -3504.111

# Row 46
$ws.Range("H46").Value = 557.3077
$ws.Range("I46").Value = 549
$ws.Range("J46").Value = 562.5
$ws.Range("K46").Value = 549
$ws.Range("L46").Value = 562.5
$ws.Range("M46").Value = -361
$ws.Range("N46").Value = -938.5

# Row 94
$ws.Range("H94").Value = 33996.668
$ws.Range("J94").Value = 33996.668
$ws.Range("L94").Value = 33996.668
$ws.Range("N94").Value = -35348.668

# Row 100
$ws.Range("H100").Value = 2130.5293
$ws.Range("I100").Value = 1212.375
$ws.Range("K100").Value = 1212.375
$ws.Range("M100").Value = -671.375

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 116666
$ws.Range("J54").Value = 116666
$ws.Range("L54").Value = 116666
$ws.Range("N54").Value = -117706

# Row 81
$ws.Range("H81").Value = 30486.666
$ws.Range("I81").Value = 1700
$ws.Range("J81").Value = 38711.43
$ws.Range("K81").Value = 3400
$ws.Range("L81").Value = 77422.86
$ws.Range("M81").Value = -2339
$ws.Range("N81").Value = -79544.86

# Row 84
$ws.Range("H84").Value = 30486.666
$ws.Range("I84").Value = 1700
$ws.Range("J84").Value = 38711.43
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 387114.3
$ws.Range("M84").Value = -11696
$ws.Range("N84").Value = -397722.3
